# Fill in the "User Story" hour estimates (column C) and the "Sprint" column (column D)
# for each task row, as part of adding user stories and sprint info to the task list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sprint = "sprint 3 "

$hours = @{
    2  = 11
    3  = 2
    4  = 2
    5  = 2
    6  = 2
    7  = 10
    8  = 8
    9  = 8
    10 = 8
    11 = 19
    12 = 9
    13 = 19
}

foreach ($row in 2..13) {
    $ws.Cells.Item($row, 3).Value = $hours[$row]
    $ws.Cells.Item($row, 4).Value = $sprint
}
